# Update SCD0013-001 "Admin SLN melakukan Modul Mapping" test-case sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab from SCD0216 -> SCD0013
$ws.Name = "SCD0013"

# TC_ID cell (B2) changes from "DGS-231" to "SCD0013-001"
$ws.Range("B2").Value = "SCD0013-001"

# Column B widens to fit the new, longer TC_ID text
$ws.Columns("B").ColumnWidth = 11.7

# Move/restore the active selection to B3
$ws.Range("B3").Select()

$wb.Save()
